$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2548.5
$ws.Range("J12").Value = 3334
$ws.Range("L12").Value = 3334
$ws.Range("N12").Value = -3674

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 62138660
$ws.Range("I86").Value = 93753064
$ws.Range("J86").Value = 7942535
$ws.Range("K86").Value = 93753064
$ws.Range("L86").Value = 7942535
$ws.Range("M86").Value = -93751941
$ws.Range("N86").Value = -7944781

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 62138660
$ws.Range("I89").Value = 93753064
$ws.Range("J89").Value = 7942535
$ws.Range("K89").Value = 468765320
$ws.Range("L89").Value = 39712675
$ws.Range("M89").Value = -468759704
$ws.Range("N89").Value = -39723907

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1451.875
$ws.Range("J92").Value = 2995
$ws.Range("L92").Value = 2995
$ws.Range("N92").Value = -5491

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 20536712
$ws.Range("I107").Value = 10228345
$ws.Range("J107").Value = 58334056
$ws.Range("K107").Value = 10228345
$ws.Range("L107").Value = 58334056
$ws.Range("M107").Value = -10226425
$ws.Range("N107").Value = -58337896

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 84800700
$ws.Range("I113").Value = 123458560
$ws.Range("J113").Value = 50008640
$ws.Range("K113").Value = 123458560
$ws.Range("L113").Value = 50008640
$ws.Range("M113").Value = -123455306
$ws.Range("N113").Value = -50015148

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 41673452
$ws.Range("J116").Value = 11803
$ws.Range("L116").Value = 11803
$ws.Range("N116").Value = -18687

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1978.8096
$ws.Range("I132").Value = 1977.45
$ws.Range("J132").Value = 2006
$ws.Range("K132").Value = 5932.35
$ws.Range("L132").Value = 6018
$ws.Range("M132").Value = -3402.35
$ws.Range("N132").Value = -11078

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1668019.5
$ws.Range("I135").Value = 2001243.8
$ws.Range("J135").Value = 1898
$ws.Range("K135").Value = 18011194.2
$ws.Range("L135").Value = 17082
$ws.Range("M135").Value = -18008659.2
$ws.Range("N135").Value = -22152

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2080.647
$ws.Range("I137").Value = 2339.2
$ws.Range("J137").Value = 1711.2858
$ws.Range("K137").Value = 7017.599999999999
$ws.Range("L137").Value = 5133.857400000001
$ws.Range("M137").Value = -4467.599999999999
$ws.Range("N137").Value = -10233.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1740008.2
$ws.Range("I32").Value = 1842297.5
$ws.Range("J32").Value = 1092.5
$ws.Range("K32").Value = 1842297.5
$ws.Range("L32").Value = 1092.5
$ws.Range("M32").Value = -1842010.5
$ws.Range("N32").Value = -1666.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7634
$ws.Range("I45").Value = 3341.4
$ws.Range("J45").Value = 12999.75
$ws.Range("K45").Value = 3341.4
$ws.Range("L45").Value = 12999.75
$ws.Range("M45").Value = -2964.4
$ws.Range("N45").Value = -13753.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3090.6736
$ws.Range("J61").Value = 7938.778
$ws.Range("L61").Value = 7938.778
$ws.Range("N61").Value = -8362.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 52381.574
$ws.Range("I74").Value = 82854
$ws.Range("J74").Value = 5500.923
$ws.Range("K74").Value = 82854
$ws.Range("L74").Value = 5500.923
$ws.Range("M74").Value = -81980
$ws.Range("N74").Value = -7248.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 52381.574
$ws.Range("I77").Value = 82854
$ws.Range("J77").Value = 5500.923
$ws.Range("K77").Value = 414270
$ws.Range("L77").Value = 27504.615
$ws.Range("M77").Value = -409902
$ws.Range("N77").Value = -36240.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3286.375
$ws.Range("I132").Value = 1992.8518
$ws.Range("K132").Value = 5978.555399999999
$ws.Range("M132").Value = -3448.555399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3090.6736
$ws.Range("J136").Value = 7938.778
$ws.Range("L136").Value = 23816.334
$ws.Range("N136").Value = -28916.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10418723
$ws.Range("J20").Value = 1774.1
$ws.Range("L20").Value = 1774.1
$ws.Range("N20").Value = -2268.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3551.1035
$ws.Range("I16").Value = 2083.7778
$ws.Range("K16").Value = 2083.7778
$ws.Range("M16").Value = -1796.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8779722
$ws.Range("I31").Value = 2731.9
$ws.Range("J31").Value = 18531934
$ws.Range("K31").Value = 2731.9
$ws.Range("L31").Value = 18531934
$ws.Range("M31").Value = -2436.9
$ws.Range("N31").Value = -18532524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8779722
$ws.Range("I34").Value = 2731.9
$ws.Range("J34").Value = 18531934
$ws.Range("K34").Value = 2731.9
$ws.Range("L34").Value = 18531934
$ws.Range("M34").Value = -2529.9
$ws.Range("N34").Value = -18532338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 848.45
$ws.Range("J94").Value = 629.0909
$ws.Range("L94").Value = 629.0909
$ws.Range("N94").Value = -1531.0909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12438.625
$ws.Range("J99").Value = 10666.333
$ws.Range("L99").Value = 10666.333
$ws.Range("N99").Value = -13662.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3551.1035
$ws.Range("I113").Value = 2083.7778
$ws.Range("K113").Value = 2083.7778
$ws.Range("M113").Value = 86.22220000000016

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12438.625
$ws.Range("J126").Value = 10666.333
$ws.Range("L126").Value = 31998.999
$ws.Range("N126").Value = -36938.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10005229
$ws.Range("I132").Value = 2248.6956
$ws.Range("J132").Value = 23538672
$ws.Range("K132").Value = 6746.0868
$ws.Range("L132").Value = 70616016
$ws.Range("M132").Value = -4216.0868
$ws.Range("N132").Value = -70621076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5763.0557
$ws.Range("I134").Value = 1364.35
$ws.Range("K134").Value = 4093.05
$ws.Range("M134").Value = -1558.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 502502500
$ws.Range("I4").Value = 505000000
$ws.Range("K4").Value = 1515000000
$ws.Range("M4").Value = -1514999888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1116.125
$ws.Range("I109").Value = 1116.125
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 3348.375
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -2308.375
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 88076.25
$ws.Range("I138").Value = 95262.37
$ws.Range("J138").Value = 9029
$ws.Range("K138").Value = 285787.11
$ws.Range("L138").Value = 27087
$ws.Range("M138").Value = -280647.11
$ws.Range("N138").Value = -37367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1648.2903
$ws.Range("J132").Value = 3277.6667
$ws.Range("L132").Value = 9833.000100000001
$ws.Range("N132").Value = -14893.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6728.0312
$ws.Range("I40").Value = 5029.2
$ws.Range("K40").Value = 5029.2
$ws.Range("M40").Value = -4893.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5180
$ws.Range("I132").Value = 3163.375
$ws.Range("K132").Value = 9490.125
$ws.Range("M132").Value = -6960.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8177.5
$ws.Range("I136").Value = 4920
$ws.Range("K136").Value = 14760
$ws.Range("M136").Value = -12210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 648.4783
$ws.Range("I107").Value = 449.53845
$ws.Range("K107").Value = 1348.61535
$ws.Range("M107").Value = 571.38465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3612.125
$ws.Range("I126").Value = 999.6
$ws.Range("J126").Value = 7966.3335
$ws.Range("K126").Value = 2998.8
$ws.Range("L126").Value = 23899.0005
$ws.Range("M126").Value = -528.8000000000002
$ws.Range("N126").Value = -28839.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14076.233
$ws.Range("I132").Value = 7716.353
$ws.Range("J132").Value = 22393
$ws.Range("K132").Value = 23149.059
$ws.Range("L132").Value = 67179
$ws.Range("M132").Value = -20619.059
$ws.Range("N132").Value = -72239

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 55030.953
$ws.Range("I136").Value = 2557.3333
$ws.Range("K136").Value = 7671.999899999999
$ws.Range("M136").Value = -5121.999899999999
